$d = $word.ActiveDocument

function Set-ParagraphXml($Index, $InnerXml) {
    $p = $d.Paragraphs.Item($Index)
    $r = $p.Range
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' + $InnerXml + '</w:p></w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($pkg)
}

# Paragraph 1: center title, keep bold run, update trailing empty run to Times/24
Set-ParagraphXml 1 (
    '<w:pPr><w:jc w:val="center"/></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="22"/><w:sz-cs w:val="22"/><w:b/></w:rPr><w:t xml:space="preserve">_Unit&#233; d&#39;enseignement_ _N_ (_UE_)</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve"></w:t></w:r>'
)

# Paragraph 2: was unjustified empty paragraph, now justified ("both") and Times/24
Set-ParagraphXml 2 (
    '<w:pPr><w:jc w:val="both"/></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve"></w:t></w:r>'
)

# Paragraph 3: right aligned, keep bold "__" run, update trailing empty run to Times/24
Set-ParagraphXml 3 (
    '<w:pPr><w:jc w:val="right"/></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="22"/><w:sz-cs w:val="22"/><w:b/></w:rPr><w:t xml:space="preserve">__</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve"></w:t></w:r>'
)

# Paragraph 4: justified empty paragraph, update run to Times/24
Set-ParagraphXml 4 (
    '<w:pPr><w:jc w:val="both"/></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve"></w:t></w:r>'
)

# Paragraph 5: justified, keep bold "__" run, update trailing empty run to Times/24
Set-ParagraphXml 5 (
    '<w:pPr><w:jc w:val="both"/></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="22"/><w:sz-cs w:val="22"/><w:b/></w:rPr><w:t xml:space="preserve">__</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve"></w:t></w:r>'
)

# Paragraph 6: justified empty paragraph, update run to Times/24
Set-ParagraphXml 6 (
    '<w:pPr><w:jc w:val="both"/></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve"></w:t></w:r>'
)

# Paragraph 7: justified, single run with "__" text, update to Times/24
Set-ParagraphXml 7 (
    '<w:pPr><w:jc w:val="both"/></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">__</w:t></w:r>'
)
